$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the location for "dal" from "cupboard over fridge" to "cupboard over sink"
$ws.Range("B9").Value = "cupboard over sink"

# Reflect the new active selection left by the edit
$ws.Range("B9").Select()
